$d = $word.ActiveDocument

# 1. Add a hanging-indent paragraph format to the (only) body paragraph
#    <w:ind w:left="1800" w:hanging="360"/>  (twips -> points: /20)
$para = $d.Paragraphs.Item(1)
$para.Range.ParagraphFormat.LeftIndent = 90
$para.Range.ParagraphFormat.FirstLineIndent = -18

# 2. Remove the "PAGE \* MERGEFORMAT" field codes from every footer of every
#    section (these produced the literal "8", "9" and "1" page-number text).
foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            for ($f = $ftr.Range.Fields.Count; $f -ge 1; $f--) {
                $ftr.Range.Fields.Item($f).Delete()
            }
        }
    }
}
